$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.509.61"
$ws.Range("E2").Value = "  -3.82%  "
$ws.Range("D3").Value = "2.909.84"
$ws.Range("E3").Value = "  -4.08%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "498.90"
$ws.Range("E5").Value = "  -2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.87"
$ws.Range("E6").Value = "  -4.84%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.419"
$ws.Range("E8").Value = "  -5.01%  "
$ws.Range("E9").Value = "  -4.50%  "
$ws.Range("E10").Value = "  -5.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.349"
$ws.Range("E11").Value = "  -4.81%  "
$ws.Range("D12").Value = "3.398.03"
$ws.Range("E12").Value = "  -4.55%  "
$ws.Range("E13").Value = "  -3.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.64"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("E15").Value = "  -5.28%  "
$ws.Range("D16").Value = "55.449.48"
$ws.Range("E16").Value = "  -4.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.00"
$ws.Range("E17").Value = "  -3.45%  "
$ws.Range("D18").Value = "2.908.77"
$ws.Range("E18").Value = "  -4.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.63"
$ws.Range("E20").Value = "  -4.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.55"
$ws.Range("E21").Value = "  -6.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.482"
$ws.Range("E23").Value = "  -3.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.26"
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").Value = "3.017.85"
$ws.Range("E25").Value = "  -4.69%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -4.88%  "
$ws.Range("D28").Value = "0.0₃0845"
$ws.Range("E28").Value = "  -9.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.34"
$ws.Range("E29").Value = "  -6.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.91"
$ws.Range("E30").Value = "  -7.15%  "
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.66"
$ws.Range("E32").Value = "  -5.66%  "
$ws.Range("E33").Value = "  -7.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.13"
$ws.Range("E34").Value = "  -4.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.36"
$ws.Range("E35").Value = "  -7.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.58"
$ws.Range("E36").Value = "  -4.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "24.54"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.18"
$ws.Range("E38").Value = "  -7.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0648"
$ws.Range("E39").Value = "  -5.32%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.19"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.68"
$ws.Range("E42").Value = "  -5.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.633"
$ws.Range("E43").Value = "  -3.79%  "
$ws.Range("D44").Value = "2.087.74"
$ws.Range("E44").Value = "  -9.26%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.32"
$ws.Range("E45").Value = "  -7.55%  "
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.87"
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.912"
$ws.Range("E47").Value = "  -7.82%  "
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.51"
$ws.Range("E49").Value = "  -5.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0833"
$ws.Range("E50").Value = "  -6.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.67"
$ws.Range("E51").Value = "  -8.72%  "
